$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'topic'
$ws.Range('B1').Value = 'variable1'
$ws.Range('C1').Value = 'variable2'
$ws.Range('D1').Value = 'variable3'
$ws.Range('E1').Value = 'question1'
$ws.Range('F1').Value = 'question2'
$ws.Range('G1').Value = 'answer1'
$ws.Range('H1').Value = 'answer2'
$ws.Range('I1').Value = 'calculation'
$ws.Range('J1').Value = 'result'
$ws.Range('K1').Value = 'response'
$ws.Range('D2').Value = 'hours_practiced'
$ws.Range('E2').Value = 'How many hours a day has Squidward practiced?'
$ws.Range('K2').Value = '"Squidward has practiced", {result}, "hours"'
$ws.Range('D3').Value = 'total'
$ws.Range('K3').Value = '{answer1}, "+", {answer2}, "is", {result}'
$ws.Range('D4').Value = 'total'
$ws.Range('K4').Value = '{answer1}, "*", {answer2}, "is", {result}'
$ws.Range('D5').Value = 'total'
$ws.Range('K5').Value = '{answer1}, "/", {answer2}, "is", {result}'
$ws.Range('D6').Value = 'total'
$ws.Range('K6').Value = '{answer1}, "-", {answer2}, "is", {result}'
$ws.Range('D7').Value = 'age_in_ten'
$ws.Range('K7').Value = '"In 10 years, you will be", {result}'
$ws.Range('D8').Value = 'age_next_year'
$ws.Range('K8').Value = '"Next year, you''ll be", {result}'
$ws.Range('D9').Value = 'age_last_year'
$ws.Range('K9').Value = '"Last year, you were", {result}'
$ws.Range('D10').Value = 'slices_left'
$ws.Range('K10').Value = '"You have", {result}, "slices left"'
$ws.Range('D11').Value = 'each_owe'
$ws.Range('K11').Value = '"You each owe £", {result}'
$ws.Range('D12').Value = 'weight_kg'
$ws.Range('K12').Value = '{answer1}, "lbs is", {result}, "kg"'
$ws.Range('D13').Value = 'fits_in'
$ws.Range('K13').Value = '"There are", {result}, {answer2},"s in", {answer1}'
$ws.Range('D14').Value = 'remainder'
$ws.Range('K14').Value = '"There are", {result}, "left over after dividing", {answer1}, "by", {answer2}'
$ws.Range('D15').Value = 'hours'
$ws.Range('K15').Value = '"There are", {result}, "hours in", {answer1}, "days"'
$ws.Range('D16').Value = 'months_to_save'
$ws.Range('K16').Value = '"It will take", {result}, "weeks to save ", {answer2}'
$ws.Range('D17').Value = 'total_spend'
$ws.Range('K17').Value = '{answer1}, "chocolate bars is", {result}'
$ws.Range('D18').Value = 'daily_spend'
$ws.Range('K18').Value = '"You will have £", {result}, "a day"'
$ws.Range('D19').Value = 'area'
$ws.Range('K19').Value = '"The area of the rectangle is", {result}'
$ws.Range('D20').Value = 'perimeter'
$ws.Range('K20').Value = '"The perimeter of the rectangle is", {result}'
$ws.Range('D21').Value = 'seconds'
$ws.Range('K21').Value = '"There are", {result}, "seconds in", {answer1}, "minutes"'
$ws.Range('D22').Value = 'average'
$ws.Range('K22').Value = '"The average is", {result}'
$ws.Range('D23').Value = 'difference'
$ws.Range('K23').Value = '"Person 1 is", {result}, "years older"'
$ws.Range('D24').Value = 'cookie_boxes'
$ws.Range('K24').Value = '{answer1}, "needs", {result}, "boxes"'
$ws.Range('D25').Value = 'total_cycled'
$ws.Range('K25').Value = '"You cycled", {result}, "km in", {answer2}, "days"'
$ws.Range('D26').Value = 'square'
$ws.Range('K26').Value = '{answer1}, "squared is", {result}'
$ws.Range('D27').Value = 'half'
$ws.Range('K27').Value = '"Half of", {answer1}, "is", {result}'
$ws.Range('D28').Value = 'double'
$ws.Range('K28').Value = '"Double", {answer1}, "is", {result}'
$ws.Range('D29').Value = 'area'
$ws.Range('K29').Value = '"The square is", {result}, "cm2"'
$ws.Range('D30').Value = 'perimeter'
$ws.Range('K30').Value = '"The square is", {result}, "cm2"'
$ws.Range('D31').Value = 'triple'
$ws.Range('K31').Value = '"Triple", {answer1}, "is", {result}'
$ws.Range('D32').Value = 'fahrenheit'
$ws.Range('K32').Value = '{answer1}, "celsius is", {result}, "fahrenheit"'
$ws.Range('B33').Value = 'age_years'
$ws.Range('D33').Value = 'age_months'
$ws.Range('K33').Value = '"You are", {result}, "months old"'
$ws.Range('D34').Value = 'area'
$ws.Range('K34').Value = '"The circle has an area of", {result}'

$ws.Range("E3").Select() | Out-Null
